$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("NewsPaper")
$ws4.Range("A28:E28").Merge()
$ws4.Range("A28").Value = "Note"
$ws4.Range("A28:E28").Font.Bold = $true
$ws4.Range("A28:E28").HorizontalAlignment = -4131
Write-Host "done"
